# Update the dSF column (F) values as per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 2
$ws.Range("F9").Value = 1
